$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "path_version"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
}

$ws.Columns.Item(10).ColumnWidth = 13.25

$ws.Range("L24").Select()
